# Added new Features to Inventory
# - Column D: per-row "Total" = Quantity * Price (PRODUCT of B:C)
# - Row 40: grand totals for Quantity, Price and Total, with the Total
#   figure highlighted in red.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row totals in column D (rows 2-10)
$ws.Range("D2").Formula = "=PRODUCT(B2:C2)"
$ws.Range("D3:D10").Formula = "=PRODUCT(B3:C3)"

# Grand-total row
$ws.Range("B40").Formula = "=SUM(B2:B10)"
$ws.Range("C40").Formula = "=SUM(C2:C10)"
$ws.Range("D40").Formula = "=SUM(D2:D10)"

# Highlight the grand-total "Total" figure in red
$ws.Range("D40").Font.Color = 255

# Match the author's final selection state
$ws.Range("B40:D40").Select()
